$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.602602124214172
$ws.Range("B1").Value = 1.636018633842468
$ws.Range("C1").Value = 1.592564344406128
$ws.Range("D1").Value = 1.917804360389709
$ws.Range("E1").Value = 2.735721588134766
